$wb = $excel.ActiveWorkbook

$dev = $wb.Worksheets.Item("2. Dev")
$wellbuilder = $wb.Worksheets.Item("3. Well Builder")

# ---------------------------------------------------------------------------
# "2. Dev" sheet (sheet3.xml)
# ---------------------------------------------------------------------------

# Rows 10-12: lowercase "x" markers become uppercase "X"
$dev.Range("A10").Value = "X"
$dev.Range("A11").Value = "X"
$dev.Range("A12").Value = "X"

# Row 13-14: text updates (no other structural change)
$dev.Range("C13").Value = "Size of chart needs to be bigger 3S-714 example plan view"
$dev.Range("C14").Value = "uncheck box for formations if none entered"

# Row 15: new text + strike-through styling on B15:C15 (item completed)
$dev.Range("C15").Value = "Casing show depths on charts?"
$dev.Range("B15:C15").Font.Strikethrough = $true

# Rows 16-18: text updates only
$dev.Range("C16").Value = "Formations names on side of the chart opposite well (if well N/S last point is - put them on the left, if + on right)"
$dev.Range("C17").Value = "close all charts when closing pad or chart view window"
$dev.Range("C18").Value = "if show is not selected do not highlight"

# Row 19: mark X, text update
$dev.Range("A19").Value = "X"
$dev.Range("C19").Value = "annotate wells"

# Row 20: text update only
$dev.Range("C20").Value = "highlight current well"

# Row 21: mark X, text update
$dev.Range("A21").Value = "X"
$dev.Range("C21").Value = "Add offset well button"

# Rows 22-26: new "X" markers in column A + text updates
$dev.Range("A22").Value = "X"
$dev.Range("C22").Value = "Plan vs Actual window"

$dev.Range("A23").Value = "X"
$dev.Range("C23").Value = "Update actual window"

$dev.Range("A24").Value = "X"
$dev.Range("C24").Value = "Update plan window"

$dev.Range("A25").Value = "X"
$dev.Range("C25").Value = "Modify importCsv.py to set as plan in databse"

$dev.Range("A26").Value = "X"
$dev.Range("C26").Value = "Update database to have planned vs actual column"

# Row 27: unchanged content ("Multilateral wells") - left as-is

# Row 28: new sub-item row
$dev.Range("A28").Value = "X"
$dev.Range("C28").Value = "3.23.1"
$dev.Range("D28").Value = "Under update directional have an add lateral button"

# Row 29: new sub-item row
$dev.Range("C29").Value = "3.23.2"
$dev.Range("D29").Value = "Add planed or actual radio putton to latter name window"

# Row 30: new sub-item row
$dev.Range("C30").Value = "3.23.3"
$dev.Range("D30").Value = "Change names in annotations to inclue lateral if not Null"

# Row 31: new top-level item (previously row 29)
$dev.Range("A31").Value = "X"
$dev.Range("B31").Value = 3.24
$dev.Range("C31").Value = "TVD vs MD plot"

# ---------------------------------------------------------------------------
# "3. Well Builder" sheet (sheet4.xml) - replaced with a new casing-builder
# outline (old "Tree Widget?" note removed)
# ---------------------------------------------------------------------------

$wellbuilder.Range("B2").Value = 1
$wellbuilder.Range("C2").Value = "Build main GUI window with buttons to add different types of casing"

$wellbuilder.Range("C3").Value = 1.1
$wellbuilder.Range("D3").Value = "Surface casing"

$wellbuilder.Range("C4").Value = 1.2
$wellbuilder.Range("D4").Value = "Intermediate Casing"

$wellbuilder.Range("C5").Value = 1.3
$wellbuilder.Range("D5").Value = "Production casing"

$wellbuilder.Range("C6").Value = 1.4
$wellbuilder.Range("D6").Value = "Liners"

$wellbuilder.Range("C7").Value = 1.5
$wellbuilder.Range("D7").Value = "Upper completions"

$wellbuilder.Range("B8").Value = 2
$wellbuilder.Range("C8").Value = "Build matplotlib window to show well diagram"

$wellbuilder.Range("B9").Value = 3
$wellbuilder.Range("C9").Value = "Build database to hold casing sizes for each well"

$wellbuilder.Range("B10").Value = 4
$wellbuilder.Range("C10").Value = "Build window that shows casing information for each section"

# ---------------------------------------------------------------------------
# Selections / active sheet - "3. Well Builder" becomes the active tab,
# "2. Dev" keeps a selection on M23 (no longer the tab-selected sheet)
# ---------------------------------------------------------------------------

$dev.Range("M23").Select() | Out-Null
$wellbuilder.Activate()
$wellbuilder.Range("D14").Select() | Out-Null
